$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Angular")

# First pass: enter the new questions in the order they were authored, so new
# shared-string entries get allocated in that same order.
$ws.Range("A26").Value = "Decoraters"
$ws.Range("A27").Value = "Directives"
$ws.Range("A28").Value = "AOT and JIT"
$ws.Range("A29").Value = "Types of Decoraters"
$ws.Range("A30").Value = "custom decoraters"
$ws.Range("A31").Value = "Dependency Injection"
$ws.Range("A32").Value = "Unit testing in angular"
$ws.Range("A33").Value = "Mocking api calls in angular for unit testing"

# Second pass: re-arrange the rows into their final order (reuses the shared
# strings just created above, so no new sharedStrings entries are added here).
$ws.Range("A26").Value = "Decoraters"
$ws.Range("A27").Value = "Types of Decoraters"
$ws.Range("A28").Value = "custom decoraters"
$ws.Range("A29").Value = "Directives"
$ws.Range("A30").Value = "Dependency Injection"
$ws.Range("A31").Value = "AOT and JIT"
$ws.Range("A32").Value = "Unit testing in angular"
$ws.Range("A33").Value = "Mocking api calls in angular for unit testing"

# Column A widened to fit the new, longer entries
$ws.Columns.Item(1).ColumnWidth = 63.43

# Make Angular the active sheet / tab, with the cursor back up at A6
$ws.Activate() | Out-Null
$ws.Range("A6").Select() | Out-Null
